$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header in P1 from "調整後利率" to "擬調利率"
$ws.Range("P1").Value = "擬調利率"

# Update the active selection to O11 (matches the recorded sheetView state)
$ws.Range("O11").Select()
